$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.298
$ws.Range("D5").Value = 0.422
$ws.Range("E5").Value = 0.457
$ws.Range("F5").Value = 0.501
$ws.Range("G5").Value = 0.5580000000000001
$ws.Range("H5").Value = 0.5620000000000001

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.298
$ws.Range("D7").Value = 0.422
$ws.Range("E7").Value = 0.457
$ws.Range("F7").Value = 0.501
$ws.Range("H7").Value = 0.5620000000000001

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.325
$ws.Range("D8").Value = 0.521
$ws.Range("E8").Value = 0.556
$ws.Range("F8").Value = 0.584
$ws.Range("G8").Value = 0.641
$ws.Range("H8").Value = 0.651

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.277
$ws.Range("C9").Value = 0.407
$ws.Range("D9").Value = 0.541
$ws.Range("E9").Value = 0.572
$ws.Range("F9").Value = 0.598
$ws.Range("G9").Value = 0.625
$ws.Range("H9").Value = 0.632
